# Updated cryptos list -- applies the refreshed Coin/Link/Price/Volume(1h)
# values captured in the source diff (row 17/18 also swap Coin+Link+Price+
# Volume, since Uniswap overtook ShibaInu in rank on this run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") often looks numeric (e.g. "1.001", "23.16"); the sheet
# stores prices as plain text, so a leading apostrophe is used (exactly as
# typing into Excel would) to keep Excel from re-typing the cell as a Number
# and silently dropping significant trailing/leading zeros.

$ws.Range("D2").Value = '29.144.82'
$ws.Range("E2").Value = '  -2.17%  '

$ws.Range("D3").Value = '1.852.31'
$ws.Range("E3").Value = '  -0.95%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'237.59"
$ws.Range("E5").Value = '  -1.47%  '

$ws.Range("D6").Value = "'0.6878"
$ws.Range("E6").Value = '  -5.38%  '

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = "'0.07707"
$ws.Range("E8").Value = '  +8.05%  '

$ws.Range("D9").Value = "'0.3035"
$ws.Range("E9").Value = '  -3.06%  '

$ws.Range("D10").Value = "'23.16"
$ws.Range("E10").Value = '  -5.01%  '

$ws.Range("D11").Value = "'0.08143"
$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").Value = '1.854.31'
$ws.Range("E12").Value = '  -1.23%  '

$ws.Range("D13").Value = "'0.7239"
$ws.Range("E13").Value = '  -2.44%  '

$ws.Range("E14").Value = '  -2.57%  '

$ws.Range("D15").Value = "'88.97"
$ws.Range("E15").Value = '  -3.72%  '

$ws.Range("D16").Value = '29.150.97'
$ws.Range("E16").Value = '  -2.15%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = "'5.734"
$ws.Range("E17").Value = '  -4.37%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'0.000007802"
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").Value = "'13.17"
$ws.Range("E19").Value = '  -1.51%  '

$ws.Range("D20").Value = "'235.18"
$ws.Range("E20").Value = '  -4.86%  '

$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").Value = '2.101.81'
$ws.Range("E22").Value = '  -1.22%  '

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").Value = "'7.554"
$ws.Range("E24").Value = '  -2.36%  '

$ws.Range("D25").Value = "'161.60"
$ws.Range("E25").Value = '  -1.25%  '

$ws.Range("D26").Value = "'8.961"
$ws.Range("E26").Value = '  -2.64%  '

$ws.Range("E27").Value = '  -7.27%  '

$ws.Range("E28").Value = '  -2.54%  '

$ws.Range("D29").Value = "'1.963"
$ws.Range("E29").Value = '  -1.90%  '

$ws.Range("D31").Value = "'4.528"
$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D33").Value = "'3.997"
$ws.Range("E33").Value = '  -4.21%  '

$ws.Range("D34").Value = "'0.05191"
$ws.Range("E34").Value = '  -2.30%  '

$ws.Range("D35").Value = "'1.177"
$ws.Range("E35").Value = '  -4.03%  '

$ws.Range("D36").Value = "'1.027"
$ws.Range("E36").Value = '  +2.93%  '

$ws.Range("D37").Value = "'0.7020"
$ws.Range("E37").Value = '  -4.77%  '

$ws.Range("D38").Value = "'2.655"
$ws.Range("E38").Value = '  -1.74%  '

$ws.Range("D40").Value = "'2.675"
$ws.Range("E40").Value = '  -2.16%  '

$ws.Range("D41").Value = "'0.9048"
$ws.Range("E41").Value = '  +2.92%  '

$ws.Range("D42").Value = '1.100.41'
$ws.Range("E42").Value = '  +5.60%  '

$ws.Range("D43").Value = "'6.001"
$ws.Range("E43").Value = '  +0.84%  '

$ws.Range("D44").Value = "'0.4269"
$ws.Range("E44").Value = '  -4.46%  '

$ws.Range("D45").Value = "'70.42"
$ws.Range("E45").Value = '  -1.07%  '

$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").Value = "'102.90"
$ws.Range("E47").Value = '  -0.85%  '

$ws.Range("D48").Value = "'1.757"
$ws.Range("E48").Value = '  -3.47%  '

$ws.Range("D49").Value = '1.997.58'
$ws.Range("E49").Value = '  -1.08%  '

$ws.Range("D50").Value = "'9.149"
$ws.Range("E50").Value = '  -4.05%  '

$ws.Range("D51").Value = "'6.924"
$ws.Range("E51").Value = '  -7.12%  '
